$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 only needs its "Tipo" column updated (Deportivo -> Deportiva); the
# rest of the row (#, Nombres, Telefono, Email) is unchanged.
$ws.Cells.Item(2, 3).Value = "Deportiva"

# Rows 3-9 are being re-ordered (grouped by Tipo: Deportiva / Publicidad /
# Television) so every column in those rows gets rewritten with the new
# values.
$newRows = @(
    @("6", "Bwin", "Deportiva", "555-3456", "support@bwin.com"),
    @("5", "Adidas", "Publicidad", "555-2345", "service@adidas.com"),
    @("8", "Puma", "Publicidad", "555-7890", "contact@puma.com"),
    @("3", "Nike", "Publicidad", "555-8765", "support@nike.com"),
    @("4", "Fox Sports", "Televisión", "555-4321", "info@foxsports.com"),
    @("2", "ESPN", "Televisión", "555-5678", "info@espn.com"),
    @("7", "Telemundo", "Televisión", "555-6543", "info@telemundo.com")
)

# Column A holds numeric-looking text ("1".."8"). Mark it as Text first so
# Excel keeps storing the new values as strings (matching the workbook's
# original shared-string typing) instead of converting them to numbers.
$ws.Range("A3:A9").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $i + 3
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# The "Tipo" column is now widest for the "Televisi\u00f3n" entries, so refresh
# its best-fit width to match the new content.
$ws.Columns.Item(3).ColumnWidth = 9.3
